$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("mslists")

# Mark "done" (column G) for accessors that now have peakList support,
# and "ionize" (column F) where applicable.
$ws.Range("F7").Value = "X"
$ws.Range("G7").Value = "X"

$ws.Range("G10").Value = "X"
$ws.Range("G11").Value = "X"
$ws.Range("G12").Value = "X"

$ws.Range("F13").Value = "X"
$ws.Range("G13").Value = "X"

$ws.Activate()
$ws.Range("G15").Select()
